# "Generate Report for Handback"
#
# The workbook tracks handback status for two files
# (0a21c7a9-...md and dc09d2d6-...md). This run removes the
# dc09d2d6 entry (row 3) from every sheet and refreshes the
# handoff/handback timestamps recorded for the remaining
# 0a21c7a9 entry on the language sheets.
#
# NOTE: this sandboxed COM engine's `Range(...).Hyperlinks.Delete()`
# (and `Hyperlinks.Item(n).Delete()`) do not scope to the range/item -
# they remove every hyperlink on the worksheet. So for each sheet we:
#   1) delete all hyperlinks,
#   2) make the data edits (timestamps, row removal),
#   3) re-create only the hyperlinks that must survive (row 2),
#      restoring the "HyperLink" look (underline + blue) by hand
#      since re-adding a hyperlink resets the cell to a fresh style.

$wb = $excel.ActiveWorkbook

$hyperlinkBlue = 15570276  # OLE BGR for #6495ED, the workbook's HyperLink font color

function Restore-HyperlinkLook($range) {
    $range.Font.Underline = $true
    $range.Font.Color = $hyperlinkBlue
}

# ---------------------------------------------------------------
# Sheet "Overview": drop row 3 (dc09d2d6 entry)
# ---------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Hyperlinks.Delete()

$wsOverview.Rows.Item(3).Delete()

$wsOverview.Hyperlinks.Add($wsOverview.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/3f9f4764cfecd0d617c278efb5b73d38cd131e83/e2e/0a21c7a9-ac5e-49a4-9d88-e7669fa35205.md", "", "", "0a21c7a9-ac5e-49a4-9d88-e7669fa35205.md")
$wsOverview.Range("A2").Value = "0a21c7a9-ac5e-49a4-9d88-e7669fa35205.md"
Restore-HyperlinkLook $wsOverview.Range("A2")

# ---------------------------------------------------------------
# Sheet "zh-cn": refresh row 2 timestamps, drop row 3
# ---------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Hyperlinks.Delete()

$wsZh.Range("E2").Value = "2016-03-18 07:18:24"
$wsZh.Range("H2").Value = "2016-03-18 07:18:41"

$wsZh.Rows.Item(3).Delete()

$wsZh.Hyperlinks.Add($wsZh.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/3f9f4764cfecd0d617c278efb5b73d38cd131e83/e2e/0a21c7a9-ac5e-49a4-9d88-e7669fa35205.md", "", "", "0a21c7a9-ac5e-49a4-9d88-e7669fa35205.md")
$wsZh.Range("A2").Value = "0a21c7a9-ac5e-49a4-9d88-e7669fa35205.md"
Restore-HyperlinkLook $wsZh.Range("A2")

$wsZh.Hyperlinks.Add($wsZh.Range("B2"), "https://github.com/OpenLocalizationTest/oltest/blob/3f9f4764cfecd0d617c278efb5b73d38cd131e83/e2e/0a21c7a9-ac5e-49a4-9d88-e7669fa35205.md", "", "", ".md")
$wsZh.Range("B2").Value = ".md"
Restore-HyperlinkLook $wsZh.Range("B2")

$wsZh.Hyperlinks.Add($wsZh.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/7337e9cfaef3f14c693aab9c0c1ea753fe2aa902/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/0a21c7a9-ac5e-49a4-9d88-e7669fa35205.6568c619aefc8b584073544e2c611af89d397d7d.zh-cn.xlf", "", "", "0a21c7a9-ac5e-49a4-9d88-e7669fa35205.6568c619aefc8b584073544e2c611af89d397d7d.zh-cn.xlf")
$wsZh.Range("D2").Value = "0a21c7a9-ac5e-49a4-9d88-e7669fa35205.6568c619aefc8b584073544e2c611af89d397d7d.zh-cn.xlf"
Restore-HyperlinkLook $wsZh.Range("D2")

$wsZh.Hyperlinks.Add($wsZh.Range("F2"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/7e41eb1152a03296fd6fbdd74564e2f934cc4dcf/e2e/0a21c7a9-ac5e-49a4-9d88-e7669fa35205.md", "", "", "0a21c7a9-ac5e-49a4-9d88-e7669fa35205.md")
$wsZh.Range("F2").Value = "0a21c7a9-ac5e-49a4-9d88-e7669fa35205.md"
Restore-HyperlinkLook $wsZh.Range("F2")

$wsZh.Hyperlinks.Add($wsZh.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/70c722538ebcddc3d7e9a384965b46c6db06f9c7/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/0a21c7a9-ac5e-49a4-9d88-e7669fa35205.6568c619aefc8b584073544e2c611af89d397d7d.zh-cn.xlf", "", "", "0a21c7a9-ac5e-49a4-9d88-e7669fa35205.6568c619aefc8b584073544e2c611af89d397d7d.zh-cn.xlf")
$wsZh.Range("G2").Value = "0a21c7a9-ac5e-49a4-9d88-e7669fa35205.6568c619aefc8b584073544e2c611af89d397d7d.zh-cn.xlf"
Restore-HyperlinkLook $wsZh.Range("G2")

# ---------------------------------------------------------------
# Sheet "de-de": refresh row 2 timestamps, drop row 3
# ---------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Hyperlinks.Delete()

$wsDe.Range("E2").Value = "2016-03-18 07:18:26"
$wsDe.Range("H2").Value = "2016-03-18 07:18:46"

$wsDe.Rows.Item(3).Delete()

$wsDe.Hyperlinks.Add($wsDe.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/3f9f4764cfecd0d617c278efb5b73d38cd131e83/e2e/0a21c7a9-ac5e-49a4-9d88-e7669fa35205.md", "", "", "0a21c7a9-ac5e-49a4-9d88-e7669fa35205.md")
$wsDe.Range("A2").Value = "0a21c7a9-ac5e-49a4-9d88-e7669fa35205.md"
Restore-HyperlinkLook $wsDe.Range("A2")

$wsDe.Hyperlinks.Add($wsDe.Range("B2"), "https://github.com/OpenLocalizationTest/oltest/blob/3f9f4764cfecd0d617c278efb5b73d38cd131e83/e2e/0a21c7a9-ac5e-49a4-9d88-e7669fa35205.md", "", "", ".md")
$wsDe.Range("B2").Value = ".md"
Restore-HyperlinkLook $wsDe.Range("B2")

$wsDe.Hyperlinks.Add($wsDe.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/0027e4c3c945d980d67524eb8dd82309ba637ead/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/0a21c7a9-ac5e-49a4-9d88-e7669fa35205.6568c619aefc8b584073544e2c611af89d397d7d.de-de.xlf", "", "", "0a21c7a9-ac5e-49a4-9d88-e7669fa35205.6568c619aefc8b584073544e2c611af89d397d7d.de-de.xlf")
$wsDe.Range("D2").Value = "0a21c7a9-ac5e-49a4-9d88-e7669fa35205.6568c619aefc8b584073544e2c611af89d397d7d.de-de.xlf"
Restore-HyperlinkLook $wsDe.Range("D2")

$wsDe.Hyperlinks.Add($wsDe.Range("F2"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/ccae32a5a8ef7264e1e4fb6706eb694e7cab819e/e2e/0a21c7a9-ac5e-49a4-9d88-e7669fa35205.md", "", "", "0a21c7a9-ac5e-49a4-9d88-e7669fa35205.md")
$wsDe.Range("F2").Value = "0a21c7a9-ac5e-49a4-9d88-e7669fa35205.md"
Restore-HyperlinkLook $wsDe.Range("F2")

$wsDe.Hyperlinks.Add($wsDe.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/a92d3434af0ba93a98a6e47caa4d83cf02a1eac0/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/ht/0a21c7a9-ac5e-49a4-9d88-e7669fa35205.6568c619aefc8b584073544e2c611af89d397d7d.de-de.xlf", "", "", "0a21c7a9-ac5e-49a4-9d88-e7669fa35205.6568c619aefc8b584073544e2c611af89d397d7d.de-de.xlf")
$wsDe.Range("G2").Value = "0a21c7a9-ac5e-49a4-9d88-e7669fa35205.6568c619aefc8b584073544e2c611af89d397d7d.de-de.xlf"
Restore-HyperlinkLook $wsDe.Range("G2")

$wsOverview.Select()
$wsOverview.Range("A1").Select()
